# Add functionality to lookup table: a separate worksheet for variable
# type (stock or flow) and its label.

$wb = $excel.ActiveWorkbook
$lookup = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Stamp an (empty, text-formatted) helper column G on "Lookup Table"
#    for the rows that correspond to entries in the new lookup sheet
#    (mirrors leftover formatting from a VLOOKUP column that was later
#    cleared out), plus H1.
# ---------------------------------------------------------------------
$gRows = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,53,55,56,57,60)
foreach ($r in $gRows) {
    $lookup.Cells.Item($r, 7).NumberFormat = "@"
}
$lookup.Range("H1").NumberFormat = "@"

# ---------------------------------------------------------------------
# 2. Add the new "Type and Label" worksheet after "Lookup Table".
# ---------------------------------------------------------------------
$typeLabel = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$typeLabel.Name = "Type and Label"

# Header row
$typeLabel.Range("A1").Value = "rec"
$typeLabel.Range("B1").Value = "type"
$typeLabel.Range("C1").Value = "label"

# Data rows: rec, type ("dollar_flow" | "flow" | "stock"), label
$rows = @(
    @("netpatrev", "dollar_flow", "net patient revenues (total revenues minus allowances & discounts)"),
    @("othinc", "dollar_flow", "other income"),
    @("opexp", "dollar_flow", "total operating expenses"),
    @("othexp", "dollar_flow", "total other expenses"),
    @("donations", "dollar_flow", "donations"),
    @("invinc", "dollar_flow", "investment income"),
    @("iphosprev", "dollar_flow", "inpatient hospital revenue"),
    @("ipgenrev", "dollar_flow", "inpatient general revenue (total of hosp, ipf, irf, snf, etc.)"),
    @("ipicrev", "dollar_flow", "inpatient intensive care type revenue (total of icu, ccu, etc.)"),
    @("iprcrev", "dollar_flow", "inpatient routine care revenue (sum of ipgenrev and ipicrev)"),
    @("ipancrev", "dollar_flow", "inpatient ancillary services revenue"),
    @("ipoprev", "dollar_flow", "inpatient outpatient services revenue"),
    @("iptotrev", "dollar_flow", "inpatient total patient revenue"),
    @("opancrev", "dollar_flow", "outpatient ancillary services revenue"),
    @("opoprev", "dollar_flow", "outpatient outpatient services revenue"),
    @("optotrev", "dollar_flow", "outpatient total patient revenues"),
    @("tottotrev", "dollar_flow", "total patient revenue (sum of iptotrev and optotrev)"),
    @("ccr", "stock", "cost to charge ratio"),
    @("totinitchcare", "dollar_flow", "total initial obligation of patients for charity care (2010 format only)"),
    @("ppaychcare", "dollar_flow", "partial payment by patients approved for charity care (2010 format only)"),
    @("nonmcbaddebt", "dollar_flow", "non-medicare & non-reimbursable medicare bad debt expense (2010 format only)"),
    @("costuccare_v2010", "dollar_flow", "cost of uncompensated care (2010 format only)"),
    @("beds_adultped", "stock", "beds - adults & peds"),
    @("availbeddays_adultped", "flow", "bed days available in rpt period"),
    @("ipbeddays_adultped", "flow", "inpatient bed days utilized"),
    @("ipdischarges_adultped", "flow", "inpatient discharges"),
    @("beds_totadultped", "stock", "beds - total adults & peds incl swing beds"),
    @("beds_total", "stock", "beds - total (inc swing + spec care beds e.g. icu, ccu, nicu)"),
    @("chguccare", "dollar_flow", "other uncompensated care charges (1996 format only)")
)

$r = 2
foreach ($row in $rows) {
    $typeLabel.Cells.Item($r, 1).Value = $row[0]
    $typeLabel.Cells.Item($r, 2).Value = $row[1]
    $typeLabel.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Column widths / formatting similar to the source sheet.
$typeLabel.Columns.Item(1).ColumnWidth = 20
$typeLabel.Range("A1:B30").NumberFormat = "@"
$typeLabel.Columns.Item(2).ColumnWidth = 10.5
$typeLabel.Columns.Item(3).ColumnWidth = 70

$lookup.Columns.Item(1).ColumnWidth = 20

# ---------------------------------------------------------------------
# 3. Restore selections / active sheet to match the final workbook
#    state.
# ---------------------------------------------------------------------
$x = $typeLabel.Range("D41").Select()
$lookup.Activate()
$y = $lookup.Range("H5").Select()
